$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 19:52"

# Countries that changed rank order (name swaps caused by updated counts below)
$ws.Range("A19").Value = "Francia"
$ws.Range("A20").Value = "Pakistan"
$ws.Range("A53").Value = "Etiopia"
$ws.Range("A54").Value = "Nigeria"
$ws.Range("A101").Value = "Maldivas"
$ws.Range("A102").Value = "Finlandia"
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 6308875
$ws.Range("C4").Value = 18138
$ws.Range("D4").Value = 3554903
$ws.Range("E4").Value = 2563571
$ws.Range("G4").Value = 437
$ws.Range("H4").Value = 190401

$ws.Range("B6").Value = 3930059
$ws.Range("C6").Value = 81091
$ws.Range("D6").Value = 3031777
$ws.Range("E6").Value = 829722
$ws.Range("G6").Value = 1074
$ws.Range("H6").Value = 68560

$ws.Range("B19").Value = 300181
$ws.Range("C19").Value = 7157
$ws.Range("D19").Value = 86963
$ws.Range("E19").Value = 182512
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 30706

$ws.Range("B20").Value = 297014
$ws.Range("C20").Value = 424
$ws.Range("D20").Value = 281925
$ws.Range("E20").Value = 8761
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 6328

$ws.Range("B23").Value = 248285
$ws.Range("C23").Value = 894
$ws.Range("E23").Value = 15790

$ws.Range("B29").Value = 123903
$ws.Range("C29").Value = 2439
$ws.Range("D29").Value = 98637
$ws.Range("E29").Value = 24281
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = 985

$ws.Range("B49").Value = 66855
$ws.Range("C49").Value = 1402
$ws.Range("D49").Value = 51223
$ws.Range("E49").Value = 14379
$ws.Range("G49").Value = 37
$ws.Range("H49").Value = 1253

$ws.Range("B53").Value = 55213
$ws.Range("C53").Value = 804
$ws.Range("D53").Value = 20283
$ws.Range("E53").Value = 34074
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 856

$ws.Range("B54").Value = 54463
$ws.Range("D54").Value = 42439
$ws.Range("E54").Value = 10997
$ws.Range("H54").Value = 1027

$ws.Range("B57").Value = 45469
$ws.Range("C57").Value = 311
$ws.Range("D57").Value = 32006
$ws.Range("E57").Value = 11934
$ws.Range("G57").Value = 6
$ws.Range("H57").Value = 1529

$ws.Range("B70").Value = 29206
$ws.Range("C70").Value = 92
$ws.Range("E70").Value = 4065

$ws.Range("B80").Value = 18963
$ws.Range("C80").Value = 588
$ws.Range("E80").Value = 13589
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 179

$ws.Range("B97").Value = 9251
$ws.Range("C97").Value = 42
$ws.Range("E97").Value = 451

$ws.Range("B101").Value = 8281
$ws.Range("C101").Value = 141
$ws.Range("D101").Value = 5483
$ws.Range("E101").Value = 2769
$ws.Range("H101").Value = 29

$ws.Range("B102").Value = 8200
$ws.Range("C102").Value = 39
$ws.Range("D102").Value = 7350
$ws.Range("E102").Value = 514
$ws.Range("H102").Value = 336

$ws.Range("B136").Value = 2533
$ws.Range("C136").Value = 1
$ws.Range("E136").Value = 1196

$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

